$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A15").Value = "LORENZA SIMONCELLI"
$ws.Range("B15").Value = "Riccardo Versini | Modium"
$ws.Range("C15").Value = "Antonio Calabrò | Avanzi"
$ws.Range("D15").Value = "Matteo Diener | U.S. Guarna"
$ws.Range("E15").Value = "Niccolò Orsi | SBARX"
$ws.Range("F15").Value = "ENRICO BORDIGNON | Pinguini Trentini"
